$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Graduate Students")

$values = @{
    "A2"  = "All students"
    "A3"  = "Male"
    "A4"  = "Female"
    "A5"  = "U.S. citizens and permanent residents"
    "A6"  = "Hispanic or Latino"
    "A7"  = "Not Hispanic or Latino"
    "A8"  = "American Indian or Alaska Native"
    "A9"  = "Asian"
    "A10" = "Black or African American"
    "A11" = "Native Hawaiian or Other Pacific Islander"
    "A12" = "White"
    "A13" = "More than one race"
    "A14" = "Unknown ethnicity and race"
    "A15" = "Temporary visa holders"
    "A16" = "Science and engineering"
    "A17" = "Science"
    "A18" = "Agricultural and veterinary sciences"
    "A19" = "Biological and biomedical sciences"
    "A20" = "Communication"
    "A21" = "Computer and information sciences"
    "A22" = "Family and consumer sciences and human sciences"
    "A23" = "Geosciences, atmospheric sciences, and ocean sciences"
    "A24" = "Mathematics and statistics"
    "A25" = "Multidisciplinary and interdisciplinary studies"
    "A26" = "Natural resources and conservation"
    "A27" = "Psychology"
    "A28" = "Physical sciences"
    "A29" = "Social sciences"
    "A30" = "Engineering"
    "A31" = "Aerospace, aeronautical, and astronautical engineering"
    "A32" = "Biological, biomedical, and biosystems engineering"
    "A33" = "Chemical, petroleum, and chemical-related engineering"
    "A34" = "Civil, environmental, transportation and related engineering fields"
    "A35" = "Electrical, electronics, communications and computer engineering"
    "A36" = "Industrial, manufacturing, systems engineering and operations research"
    "A37" = "Mechanical engineering"
    "A38" = "Metallurgical, mining, materials and related engineering fields"
    "A39" = "Other engineering"
    "A40" = "Health"
    "A41" = "Clinical medicine"
    "A42" = "Other health"
}

foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}
